$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.579.85"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "3.456.58"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "575.02"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "161.54"
$ws.Range("E6").Value = "  +3.68%  "

$ws.Range("D8").Value = "3.457.37"
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  +8.42%  "

$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("E11").Value = "  +3.85%  "

$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "4.053.97"
$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("E14").Value = "  -2.90%  "

$ws.Range("E15").Value = "  +5.04%  "

$ws.Range("D16").Value = "28.76"
$ws.Range("E16").Value = "  +6.21%  "

$ws.Range("D17").Value = "64.564.53"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "3.450.65"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  +2.07%  "

$ws.Range("D21").Value = "389.76"
$ws.Range("E21").Value = "  +0.69%  "

$ws.Range("E22").Value = "  -3.47%  "

$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("D24").Value = "72.86"
$ws.Range("E24").Value = "  +2.73%  "

$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  +19.25%  "

$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "6.20"
$ws.Range("E30").Value = "  +9.93%  "

$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  +6.05%  "

$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").Value = "6.54"
$ws.Range("E33").Value = "  -0.87%  "

$ws.Range("D34").Value = "23.58"
$ws.Range("E34").Value = "  +1.95%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").Value = "7.08"
$ws.Range("E36").Value = "  +5.89%  "

$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").Value = "161.50"
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("D40").Value = "0.0771"
$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("D41").Value = "27.44"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").Value = "2.926.12"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("E43").Value = "  +5.86%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "42.77"
$ws.Range("E44").Value = "  +3.21%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0317"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("D47").Value = "24.09"
$ws.Range("E47").Value = "  +7.30%  "

$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").Value = "0.874"
$ws.Range("E49").Value = "  +7.32%  "

$ws.Range("D50").Value = "2.18"
$ws.Range("E50").Value = "  +12.14%  "

$ws.Range("E51").Value = "  +3.27%  "
